$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'268.43"
$ws.Range("D3").Value = "'22.66"
$ws.Range("D4").Value = "'6.330"
$ws.Range("D5").Value = "'0.06179"
$ws.Range("D6").Value = "'3.661"
$ws.Range("D7").Value = "'6.665"
$ws.Range("D8").Value = "'1.386"
$ws.Range("D9").Value = "'0.8318"
$ws.Range("D10").Value = "'0.01364"
$ws.Range("D11").Value = "'0.1598"
$ws.Range("D12").Value = "'0.08295"
$ws.Range("D13").Value = "'0.03544"
$ws.Range("D14").Value = "'0.03260"
$ws.Range("D15").Value = "'4.073"
$ws.Range("D16").Value = "'0.09304"
$ws.Range("D17").Value = "'0.001635"
$ws.Range("D18").Value = "'0.04760"
$ws.Range("D19").Value = "'0.006369"
$ws.Range("D20").Value = "'0.005655"
$ws.Range("D21").Value = "'0.001078"
$ws.Range("D22").Value = "'0.0001503"
$ws.Range("D23").Value = "'3.727"
$ws.Range("D24").Value = "'2.324"
$ws.Range("D26").Value = "'0.1234"
$ws.Range("D27").Value = "'0.0002710"
$ws.Range("D40").Value = "'0.04719"
$ws.Range("D41").Value = "'0.006979"
$ws.Range("D42").Value = "'0.1158"
$ws.Range("D43").Value = "'0.003524"
$ws.Range("D44").Value = "'0.01181"
$ws.Range("D45").Value = "'0.00006267"
$ws.Range("D46").Value = "'0.0009921"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D48").Value = "'0.7838"
$ws.Range("D49").Value = "'0.002316"
$ws.Range("D50").Value = "'0.00002405"
$ws.Range("D51").Value = "'0.01243"
